$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh footer timestamp (07:42 -> 08:12)
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 08:12"

# Apply per-cell updates: newly inserted countries shift several rows
# (country name swaps) and refreshed case/death counts for the affected rows.
# Row 4
$ws.Cells.Item(4, 2).Value = 85612
$ws.Cells.Item(4, 3).Value = 177
$ws.Cells.Item(4, 5).Value = 82443
# Row 8
$ws.Cells.Item(8, 2).Value = 47278
$ws.Cells.Item(8, 3).Value = 3340
$ws.Cells.Item(8, 5).Value = 41324
$ws.Cells.Item(8, 7).Value = 14
$ws.Cells.Item(8, 8).Value = 281
# Row 36
$ws.Cells.Item(36, 5).Value = 911
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 24
# Row 39
$ws.Cells.Item(39, 5).Value = 913
$ws.Cells.Item(39, 7).Value = 2
$ws.Cells.Item(39, 8).Value = 2
# Row 41
$ws.Cells.Item(41, 4).Value = 42
$ws.Cells.Item(41, 5).Value = 823
# Row 44
$ws.Cells.Item(44, 2).Value = 753
$ws.Cells.Item(44, 3).Value = 26
$ws.Cells.Item(44, 4).Value = 67
$ws.Cells.Item(44, 5).Value = 666
# Row 143
$ws.Cells.Item(143, 1).Value = "Uganda"
# Row 144
$ws.Cells.Item(144, 1).Value = "Nueva Caledonia"
# Row 152
$ws.Cells.Item(152, 1).Value = "Dominica"
# Row 153
$ws.Cells.Item(153, 1).Value = "San Martin (Parte Francesa)"
# Row 156
$ws.Cells.Item(156, 1).Value = "Groenlandia"
$ws.Cells.Item(156, 2).Value = 9
$ws.Cells.Item(156, 3).Value = 3
$ws.Cells.Item(156, 4).Value = 2
$ws.Cells.Item(156, 5).Value = 7
# Row 157
$ws.Cells.Item(157, 1).Value = "Surinam"
# Row 158
$ws.Cells.Item(158, 1).Value = "Haiti"
$ws.Cells.Item(158, 5).Value = 8
$ws.Cells.Item(158, 8).Value = 0
# Row 159
$ws.Cells.Item(159, 1).Value = "Islas Caimanes"
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 7
$ws.Cells.Item(159, 8).Value = 1
# Row 160
$ws.Cells.Item(160, 1).Value = "Namibia"
$ws.Cells.Item(160, 2).Value = 8
$ws.Cells.Item(160, 4).Value = 2
$ws.Cells.Item(160, 5).Value = 6
# Row 162
$ws.Cells.Item(162, 1).Value = "Mozambique"
# Row 163
$ws.Cells.Item(163, 1).Value = "Antigua y Barbuda"
# Row 164
$ws.Cells.Item(164, 1).Value = "Granada"
$ws.Cells.Item(164, 5).Value = 7
$ws.Cells.Item(164, 8).Value = 0
# Row 165
$ws.Cells.Item(165, 1).Value = "Gabon"
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 5).Value = 6
# Row 166
$ws.Cells.Item(166, 1).Value = "Curazao"
$ws.Cells.Item(166, 2).Value = 7
$ws.Cells.Item(166, 4).Value = 2
$ws.Cells.Item(166, 5).Value = 4
$ws.Cells.Item(166, 8).Value = 1
# Row 167
$ws.Cells.Item(167, 1).Value = "Eritrea"
# Row 168
$ws.Cells.Item(168, 1).Value = "Suazilandia"
# Row 169
$ws.Cells.Item(169, 1).Value = "Benin"
# Row 170
$ws.Cells.Item(170, 1).Value = "Laos"
$ws.Cells.Item(170, 4).Value = 0
$ws.Cells.Item(170, 5).Value = 6
# Row 172
$ws.Cells.Item(172, 1).Value = "Montserrat"
# Row 173
$ws.Cells.Item(173, 1).Value = "Siria"
# Row 174
$ws.Cells.Item(174, 1).Value = "Fiyi"
# Row 175
$ws.Cells.Item(175, 1).Value = "Cabo Verde"
# Row 176
$ws.Cells.Item(176, 1).Value = "Guyana"
# Row 177
$ws.Cells.Item(177, 1).Value = "Santa Sede"
# Row 178
$ws.Cells.Item(178, 1).Value = "Mali"
# Row 179
$ws.Cells.Item(179, 1).Value = "Congo"
# Row 180
$ws.Cells.Item(180, 1).Value = "Guinea"
# Row 181
$ws.Cells.Item(181, 1).Value = "Angola"
# Row 183
$ws.Cells.Item(183, 1).Value = "Mauritania"
$ws.Cells.Item(183, 3).Value = 0
# Row 184
$ws.Cells.Item(184, 1).Value = "Republica de Africa Central"
# Row 185
$ws.Cells.Item(185, 1).Value = "San Martin (Parte Holandesa)"
# Row 186
$ws.Cells.Item(186, 1).Value = "San Bartolome"
# Row 187
$ws.Cells.Item(187, 1).Value = "Republica del Chad"
# Row 188
$ws.Cells.Item(188, 1).Value = "Butan"
$ws.Cells.Item(188, 3).Value = 1
# Row 189
$ws.Cells.Item(189, 1).Value = "Nepal"
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 8).Value = 0
# Row 190
$ws.Cells.Item(190, 1).Value = "Sudan"
# Row 191
$ws.Cells.Item(191, 1).Value = "Santa Lucia"
$ws.Cells.Item(191, 4).Value = 1
$ws.Cells.Item(191, 8).Value = 0
# Row 192
$ws.Cells.Item(192, 1).Value = "Gambia"
$ws.Cells.Item(192, 4).Value = 0
$ws.Cells.Item(192, 8).Value = 1
# Row 193
$ws.Cells.Item(193, 1).Value = "Zimbabue"
$ws.Cells.Item(193, 4).Value = 0
$ws.Cells.Item(193, 8).Value = 1
# Row 194
$ws.Cells.Item(194, 1).Value = "Belice"
# Row 196
$ws.Cells.Item(196, 1).Value = "Guinea-Bisau"
# Row 197
$ws.Cells.Item(197, 1).Value = "San Cristobal y Nieves"
# Row 198
$ws.Cells.Item(198, 1).Value = "Islas Turcas y Caicos"
# Row 199
$ws.Cells.Item(199, 1).Value = "Islas Virgenes Britanicas"
# Row 200
$ws.Cells.Item(200, 1).Value = "Anguila"
# Row 202
$ws.Cells.Item(202, 1).Value = "San Vicente y las Granadinas"
# Row 204
$ws.Cells.Item(204, 1).Value = "Libia"
